# Fruta / hortaliza, semanal
# Insert a new weekly price record for Femacal de La Calera - Frutilla.
# This pushes the existing rows 127-170 down to 128-171 and populates the
# newly inserted row 127 with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at sheet row 127 (shifts rows 127:170 -> 128:171)
$ws.Rows.Item(127).Insert()

# Populate the new row with the new record
$ws.Cells.Item(127, 1).Value = 3
$ws.Cells.Item(127, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(127, 3).Value = 'Coquimbo'
$ws.Cells.Item(127, 4).Value = '11/24/2021'
$ws.Cells.Item(127, 5).Value = 5
$ws.Cells.Item(127, 6).Value = 'Fruta'
$ws.Cells.Item(127, 7).Value = 100101
$ws.Cells.Item(127, 8).Value = 'Berries'
$ws.Cells.Item(127, 9).Value = 100112025
$ws.Cells.Item(127, 10).Value = 'Frutilla'
$ws.Cells.Item(127, 11).Value = 'Sin especificar'
$ws.Cells.Item(127, 12).Value = 'Primera'
$ws.Cells.Item(127, 13).Value = 480
$ws.Cells.Item(127, 14).Value = 5500
$ws.Cells.Item(127, 15).Value = 6000
$ws.Cells.Item(127, 16).Value = 5760
$ws.Cells.Item(127, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(127, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(127, 19).Value = 823
$ws.Cells.Item(127, 20).Value = 7
